$PT_OBJ_NEW = "Dar ao futuro engenheiro os conceitos fundamentais relacionados ao escoamento de fluidos e desenvolver as equações de conservação de massa, energia e quantidade de movimento. Os conceitos e modelos matemáticos estudados servem de base para a compreensão dos processos produtivos que envolvam a transferência de fluidos e para as disciplinas de Operações Unitárias que estudam os princípios destas operações."

$PT_SHORT_SYLLABUS_NEW = "1) Bases conceituais para o estudo dos Fenômenos de transporte `n2) Propriedades gerais dos fluidos `n3) Cinemática dos fluidos:. `n4) Equações de Conservação na forma Integral: `n5) Equações Diferenciais do Escoamento de Fluidos: `n6) Teoria da Camada Limite: `n7) Escoamento em tubos:"

$PT_FULL_SYLLABUS_NEW = "1) Bases conceituais para o estudo dos Fenômenos de transporte `nFluidos e a hipótese do contínuo. Importância da análise dimensional e uso dos números adimensionais. Leis básicas para transferência de massa, calor e quantidade de movimento. Lei geral para os fenômenos de transporte. Difusividade molecular, condutividade térmica e viscosidade. Transporte simultâneo de massa, calor e quantidade de movimento. Formulação integral e diferencial. `n2) Propriedades gerais dos fluidos: Massa específica, peso específico, volume específico. Tensão e Pressão. Fluidos Newtonianos e não Newtonianos. Viscosidade. Tensão superficial e Capilaridade. Módulo de elasticidade volumétrica e compressibilidade. `n3) Cinemática dos fluidos: Descrição de um Fluido em Movimento: Método de Euler e Lagrange - Campo de escoamento de um fluido - Escoamento permanente e transiente - Trajetórias e Linhas de corrente - Sistema e volume de controle - Escoamentos unidimensionais e bidimensionais. Escoamento uniforme. Escoamento laminar e turbulento: N° de Reynolds. `n4) Equações de Conservação na forma Integral: Fluxo de uma grandeza. Conservação da Massa, continuidade. Formas específicas para a expressão integral. Conservação da quantidade de movimento linear. Conservação da Energia. Equação de Bernoulli. Aplicações `n5) Equações Diferenciais do Escoamento de Fluidos: Equação da conservação da massa e continuidade. Equação da energia. Equação de Navier-Stokes. Aplicações `n6) Teoria da Camada Limite: Definição de camada limite . Camada limite laminar e turbulenta. Camada limite sobre uma placa plana. Aplicações `n7) Escoamento em tubos: Escoamento Laminar e turbulento. Coeficiente de atrito. Região turbulenta e de transição. Diagramas de Moody e Von Karman . Equação da energia com equipamentos de transporte. Perda de carga em acidentes. Diâmetro equivalente."

$BIBLIOGRAFIA_NEW = "1) YONG, D. F.; OKIISHI, T. H.; MUNSON, B.R. Fundamentos da mecânica dos fluidos. São Paulo: Edgard Blucher `n2) BRUNETTI, F. Mecânica dos fluídos. São Paulo: Pearson Education. `n3) FOX, Robert W. Introdução à mecânica dos fluídos. Rio de Janeiro: LTC. `n4) WHITE, Frank M. Mecânica dos fluídos. Rio de Janeiro: Mcgraw-hill Interamericana. `n5) COULSON, J. M.; RICHARDSON, J.F. Chemical engineering . Oxford: Butterworth Heinemann. Volume 1: Fluid Flow, Heat Transfer and Mass Transfer"


$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13. This shifts the old rows 13-24 down to 14-25,
# giving each "label" row (e.g. "Programa resumido:", "Método:", etc.) its own
# freshly-placed value row immediately below/within, matching the corrected layout.
$ws.Rows.Item(13).Insert()

# Row 10 (Objetivos:) - replace the (misplaced) value with the real PT objectives text.
$ws.Range("B10").Value = $PT_OBJ_NEW
$ws.Range("C10").Value = $PT_OBJ_NEW

# Row 13 (new, blank after insert) - now holds "Docentes responsáveis:" value (moved up from old row 10/col).
$ws.Range("B13").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C13").Value = "5816812 - João Paulo Alves Silva"

# Row 14 (was old row 13: "Programa resumido:" / "Semestral") - replace value with the new PT short syllabus text.
$ws.Range("B14").Value = $PT_SHORT_SYLLABUS_NEW
$ws.Range("C14").Value = $PT_SHORT_SYLLABUS_NEW

# Row 16 (was old row 15: "Programa:" / stray date) - replace value with the new PT full syllabus text.
$ws.Range("B16").Value = $PT_FULL_SYLLABUS_NEW
$ws.Range("C16").Value = $PT_FULL_SYLLABUS_NEW

# Row 19 (was old row 18: "Método:" / misplaced JP Silva text) - replace with the exam-application text.
$ws.Range("B19").Value = "Aplicação de 2 provas, P1 e P2."
$ws.Range("C19").Value = "Aplicação de 2 provas, P1 e P2."

# Row 20 (was old row 19: "Critério:" / exam-application text) - replace with the "média do período" text.
$ws.Range("B20").Value = "A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)"
$ws.Range("C20").Value = "A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)"

# Row 22 (was old row 21: "Bibliografia:" / "média final após a recuperação" text) - replace with the real bibliography text.
$ws.Range("B22").Value = $BIBLIOGRAFIA_NEW
$ws.Range("C22").Value = $BIBLIOGRAFIA_NEW
